$wb = $excel.ActiveWorkbook

function Set-Cell($ws, $ref, $val) {
    if ($null -eq $val) {
        $ws.Range($ref).ClearContents()
    } else {
        $ws.Range($ref).Value = $val
    }
}

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 68
Set-Cell $ws "H68" 0
Set-Cell $ws "J68" 0
Set-Cell $ws "L68" 0
Set-Cell $ws "N68" $null
# Row 71
Set-Cell $ws "H71" 0
Set-Cell $ws "J71" 0
Set-Cell $ws "L71" 0
Set-Cell $ws "N71" $null
# Row 76
Set-Cell $ws "H76" 3228.5715
Set-Cell $ws "I76" 3227.7778
Set-Cell $ws "K76" 3227.7778
Set-Cell $ws "M76" -2912.7778
# Row 79
Set-Cell $ws "H79" 3228.5715
Set-Cell $ws "I79" 3227.7778
Set-Cell $ws "K79" 3227.7778
Set-Cell $ws "M79" -2135.7778
# Row 135
Set-Cell $ws "H135" 19232622
Set-Cell $ws "I135" 26317018
Set-Cell $ws "J135" 3543.4285
Set-Cell $ws "K135" 236853162
Set-Cell $ws "L135" 31890.8565
Set-Cell $ws "M135" -236850627
Set-Cell $ws "N135" -36960.8565
# Row 137
Set-Cell $ws "H137" 2875540.8
Set-Cell $ws "I137" 3969805.2
Set-Cell $ws "J137" 3097
Set-Cell $ws "K137" 11909415.6
Set-Cell $ws "L137" 9291
Set-Cell $ws "M137" -11906865.6
Set-Cell $ws "N137" -14391
# Row 138
Set-Cell $ws "H138" 4445.9546
Set-Cell $ws "I138" 4961.9165
Set-Cell $ws "J138" 4364.487
Set-Cell $ws "K138" 14885.7495
Set-Cell $ws "L138" 13093.461
Set-Cell $ws "M138" -9745.749500000002
Set-Cell $ws "N138" -23373.461
# Row 141
Set-Cell $ws "H141" 3204.2666
Set-Cell $ws "I141" 1170.1818
Set-Cell $ws "J141" 8798
Set-Cell $ws "K141" 3510.5454
Set-Cell $ws "L141" 26394
Set-Cell $ws "M141" 1669.4546
Set-Cell $ws "N141" -36754

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
Set-Cell $ws "H32" 23278058
Set-Cell $ws "I32" 27048176
Set-Cell $ws "J32" 28999.834
Set-Cell $ws "K32" 27048176
Set-Cell $ws "L32" 28999.834
Set-Cell $ws "M32" -27047889
Set-Cell $ws "N32" -29573.834

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
Set-Cell $ws "H94" 67121.87
Set-Cell $ws "I94" 100371.8
Set-Cell $ws "J94" 622
Set-Cell $ws "K94" 100371.8
Set-Cell $ws "L94" 622
Set-Cell $ws "M94" -99920.8
Set-Cell $ws "N94" -1524
# Row 134
Set-Cell $ws "H134" 3818.2285
Set-Cell $ws "I134" 3876.3572
Set-Cell $ws "J134" 3585.7144
Set-Cell $ws "K134" 11629.0716
Set-Cell $ws "L134" 10757.1432
Set-Cell $ws "M134" -9094.0716
Set-Cell $ws "N134" -15827.1432

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
Set-Cell $ws "H31" 7344.977
Set-Cell $ws "I31" 1108.0952
Set-Cell $ws "J31" 9329.439
Set-Cell $ws "K31" 1108.0952
Set-Cell $ws "L31" 9329.439
Set-Cell $ws "M31" -813.0952
Set-Cell $ws "N31" -9919.439
# Row 34
Set-Cell $ws "H34" 7344.977
Set-Cell $ws "I34" 1108.0952
Set-Cell $ws "J34" 9329.439
Set-Cell $ws "K34" 1108.0952
Set-Cell $ws "L34" 9329.439
Set-Cell $ws "M34" -906.0952
Set-Cell $ws "N34" -9733.439
# Row 96
Set-Cell $ws "H96" 49312
Set-Cell $ws "J96" 49312
Set-Cell $ws "L96" 49312
Set-Cell $ws "N96" -54804
# Row 107
Set-Cell $ws "H107" 31250750
Set-Cell $ws "I107" 62500000
Set-Cell $ws "J107" 1500
Set-Cell $ws "K107" 62500000
Set-Cell $ws "L107" 1500
Set-Cell $ws "M107" -62498080
Set-Cell $ws "N107" -5340
# Row 132
Set-Cell $ws "H132" 44446652
Set-Cell $ws "I132" 55557570
Set-Cell $ws "J132" 27780278
Set-Cell $ws "K132" 166672710
Set-Cell $ws "L132" 83340834
Set-Cell $ws "M132" -166670180
Set-Cell $ws "N132" -83345894
# Row 134
Set-Cell $ws "H134" 3626041.2
Set-Cell $ws "I134" 4169595
Set-Cell $ws "J134" 2349.111
Set-Cell $ws "K134" 12508785
Set-Cell $ws "L134" 7047.333
Set-Cell $ws "M134" -12506250
Set-Cell $ws "N134" -12117.333

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
Set-Cell $ws "H5" 1200.8306
Set-Cell $ws "I5" 920.59186
Set-Cell $ws "K5" 2761.77558
Set-Cell $ws "M5" -2649.77558
# Row 92
Set-Cell $ws "H92" 290
Set-Cell $ws "I92" 289.66666
Set-Cell $ws "K92" 868.9999799999999
Set-Cell $ws "M92" 379.0000200000001
# Row 113
Set-Cell $ws "H113" 620.6269
Set-Cell $ws "I113" 599.5227
Set-Cell $ws "J113" 661
Set-Cell $ws "K113" 1798.5681
Set-Cell $ws "L113" 1983
Set-Cell $ws "M113" 371.4319
Set-Cell $ws "N113" -6323
# Row 122
Set-Cell $ws "H122" 2692.5762
Set-Cell $ws "I122" 465.5
Set-Cell $ws "J122" 4447.242
Set-Cell $ws "K122" 4189.5
Set-Cell $ws "L122" 40025.178
Set-Cell $ws "M122" -1739.5
Set-Cell $ws "N122" -44925.178
# Row 131
Set-Cell $ws "H131" 2902.4092
Set-Cell $ws "J131" 3145.0679
Set-Cell $ws "L131" 9435.2037
Set-Cell $ws "N131" -19515.2037
# Row 135
Set-Cell $ws "H135" 1200.8306
Set-Cell $ws "I135" 920.59186
Set-Cell $ws "K135" 8285.32674
Set-Cell $ws "M135" -5750.32674
# Row 137
Set-Cell $ws "H137" 11948176
Set-Cell $ws "I137" 12149
Set-Cell $ws "J137" 19671486
Set-Cell $ws "K137" 36447
Set-Cell $ws "L137" 59014458
Set-Cell $ws "M137" -31347
Set-Cell $ws "N137" -59024658

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 39
Set-Cell $ws "H39" 58261
Set-Cell $ws "J39" 58261
Set-Cell $ws "L39" 58261
Set-Cell $ws "N39" -59325
# Row 104
Set-Cell $ws "H104" 33000
Set-Cell $ws "J104" 33000
Set-Cell $ws "L104" 33000
Set-Cell $ws "N104" -39988
# Row 107
Set-Cell $ws "H107" 381.15384
Set-Cell $ws "I107" 346.25
Set-Cell $ws "J107" 800
Set-Cell $ws "K107" 346.25
Set-Cell $ws "L107" 800
Set-Cell $ws "M107" 1573.75
Set-Cell $ws "N107" -4640
# Row 122
Set-Cell $ws "H122" 1963.5
Set-Cell $ws "I122" 1896.5454
Set-Cell $ws "J122" 2700
Set-Cell $ws "K122" 5689.6362
Set-Cell $ws "L122" 8100
Set-Cell $ws "M122" -3239.6362
Set-Cell $ws "N122" -13000
# Row 126
Set-Cell $ws "H126" 3836.842
Set-Cell $ws "J126" 3990.7273
Set-Cell $ws "L126" 11972.1819
Set-Cell $ws "N126" -16912.1819
# Row 132
Set-Cell $ws "H132" 26321114
Set-Cell $ws "I132" 41673012
Set-Cell $ws "J132" 3574.1428
Set-Cell $ws "K132" 125019036
Set-Cell $ws "L132" 10722.4284
Set-Cell $ws "M132" -125016506
Set-Cell $ws "N132" -15782.4284

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 94
Set-Cell $ws "H94" 77076.664
Set-Cell $ws "J94" 77076.664
Set-Cell $ws "L94" 77076.664
Set-Cell $ws "N94" -78428.664
# Row 132
Set-Cell $ws "H132" 3394.7026
Set-Cell $ws "I132" 2818.4375
Set-Cell $ws "J132" 3833.762
Set-Cell $ws "K132" 8455.3125
Set-Cell $ws "L132" 11501.286
Set-Cell $ws "M132" -5925.3125
Set-Cell $ws "N132" -16561.286

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
Set-Cell $ws "H132" 6947005.5
Set-Cell $ws "I132" 2660.7307
Set-Cell $ws "J132" 18231566
Set-Cell $ws "K132" 7982.1921
Set-Cell $ws "L132" 54694698
Set-Cell $ws "M132" -5452.1921
Set-Cell $ws "N132" -54699758
# Row 136
Set-Cell $ws "H136" 3922.5264
Set-Cell $ws "I136" 4409.88
Set-Cell $ws "J136" 2985.3076
Set-Cell $ws "K136" 13229.64
Set-Cell $ws "L136" 8955.9228
Set-Cell $ws "M136" -10679.64
Set-Cell $ws "N136" -14055.9228
